$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-11, columns A-L
$data = @(
    @(29, "Male",   18, 9,  0,  30, "Premium",  "Quarterly", 930, 18, $true, 0),
    @(39, "Male",   12, 5,  7,  4,  "Standard", "Quarterly", 969, 13, $true, 1),
    @(30, "Female", 39, 14, 5,  18, "Standard", "Annual",    932, 17, $true, 2),
    @(23, "Male",   32, 20, 5,  8,  "Basic",    "Monthly",   617, 20, $true, 2),
    @(65, "Female", 49, 1,  10, 8,  "Basic",    "Monthly",   557, 6,  $true, 3),
    @(55, "Female", 14, 4,  6,  18, "Basic",    "Quarterly", 185, 3,  $true, 3),
    @(58, "Male",   38, 21, 7,  7,  "Standard", "Monthly",   396, 29, $true, 3),
    @(58, "Female", 49, 12, 3,  16, "Standard", "Quarterly", 821, 24, $true, 3),
    @(55, "Female", 37, 8,  4,  15, "Premium",  "Annual",    445, 30, $true, 3),
    @(64, "Female", 3,  25, 2,  11, "Standard", "Quarterly", 415, 29, $true, 3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
